# Generate Report for Handback
#
# Updates the localization-status report after a handback: the overview
# status moves from "Ready for handoff" to "Handed back: in sync with
# en-US", the per-language sheets gain "Latest Target File" / "Latest
# Handback File" links + a "Latest Handback DateTime", and the columns
# that now hold longer hyperlink text are widened.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# Character widths are stored internally as (rounded-character-width + 5/6).
# These inputs land exactly on 40 characters (9/10) or as close as the
# model allows to ~29.98 characters (3/5/6) once stored.
$wideColWidth = 29.166666666666668
$fullColWidth = 39.166666666666664

# ---------------------------------------------------------------------
# Overview sheet: status text for both rows + widen the status columns
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $handedBack
$overview.Range("F2").Value = $handedBack
$overview.Range("E3").Value = $handedBack
$overview.Range("F3").Value = $handedBack

$overview.Columns.Item(5).ColumnWidth = $wideColWidth
$overview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("C3").Value = $handedBack

$zhcn.Columns.Item(3).ColumnWidth = $wideColWidth
$zhcn.Columns.Item(9).ColumnWidth = $fullColWidth
$zhcn.Columns.Item(10).ColumnWidth = $fullColWidth

# Rebuild the hyperlinks in display order (A2, I2, A3, I3) so the new
# "Latest Target File" links land between the existing "Source File
# Name" links, matching how Excel renumbers hyperlink relationship ids
# on save.
$zhcn.Hyperlinks.Delete()

$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md", "", "", "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md", "", "", "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/a7042400-8349-4ef4-9d00-803d03e8f618.md", "", "", "a7042400-8349-4ef4-9d00-803d03e8f618.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/a7042400-8349-4ef4-9d00-803d03e8f618.md", "", "", "a7042400-8349-4ef4-9d00-803d03e8f618.md") | Out-Null

# Restore the "Source File Name" hyperlink look (A2/A3 already had it;
# give the new "Latest Target File" cells the same look).
$zhcn.Range("A2").Font.Underline = 2
$zhcn.Range("A2").Font.Color = 15570276
$zhcn.Range("A3").Font.Underline = 2
$zhcn.Range("A3").Font.Color = 15570276
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = 15570276

# "Latest Handback File" (xlf produced on handback)
$zhcn.Range("J2").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.684d8844e0884ae608929bad0eabacf861d159b2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.ac66859d14d4501f31d5e816897883c12bcfcd2e.zh-cn.xlf"

# "Latest Handback DateTime" - same handback run for both rows
$zhcn.Range("K2").Value = "2016-08-15 20:44:29"
$zhcn.Range("K3").Value = "2016-08-15 20:44:29"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede.Range("C2").Value = $handedBack
$dede.Range("C3").Value = $handedBack

$dede.Columns.Item(3).ColumnWidth = $wideColWidth
$dede.Columns.Item(9).ColumnWidth = $fullColWidth
$dede.Columns.Item(10).ColumnWidth = $fullColWidth

$dede.Hyperlinks.Delete()

$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md", "", "", "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md", "", "", "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/a7042400-8349-4ef4-9d00-803d03e8f618.md", "", "", "a7042400-8349-4ef4-9d00-803d03e8f618.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b04e7a825218e3d74ff95af77cea65a52aa8612/e2e/a7042400-8349-4ef4-9d00-803d03e8f618.md", "", "", "a7042400-8349-4ef4-9d00-803d03e8f618.md") | Out-Null

$dede.Range("A2").Font.Underline = 2
$dede.Range("A2").Font.Color = 15570276
$dede.Range("A3").Font.Underline = 2
$dede.Range("A3").Font.Color = 15570276
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = 15570276

$dede.Range("J2").Value = "6cc6d73e-a009-46b1-b8fd-922c4b59b5c0.684d8844e0884ae608929bad0eabacf861d159b2.de-de.xlf"
$dede.Range("J3").Value = "a7042400-8349-4ef4-9d00-803d03e8f618.ac66859d14d4501f31d5e816897883c12bcfcd2e.de-de.xlf"

# de-de's handback run finished later than zh-cn's
$dede.Range("K2").Value = "2016-08-15 20:44:38"
$dede.Range("K3").Value = "2016-08-15 20:44:38"
